$wb = $excel.ActiveWorkbook

# Sheet1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 521
$ws1.Cells.Item(4, 6).Value = 393
$ws1.Cells.Item(5, 6).Value = 497
$ws1.Cells.Item(6, 6).Value = 935
$ws1.Cells.Item(8, 6).Value = 12
$ws1.Cells.Item(9, 6).Value = 975
$ws1.Cells.Item(10, 6).Value = 767
$ws1.Cells.Item(11, 6).Value = 212
$ws1.Cells.Item(14, 6).Value = 797
$ws1.Cells.Item(15, 6).Value = 262
$ws1.Cells.Item(16, 6).Value = 565
$ws1.Cells.Item(17, 6).Value = 496
$ws1.Cells.Item(20, 6).Value = 438
$ws1.Cells.Item(21, 6).Value = 1134
$ws1.Cells.Item(22, 6).Value = 2823
$ws1.Cells.Item(23, 6).Value = 1342
$ws1.Cells.Item(24, 6).Value = 670
$ws1.Cells.Item(25, 6).Value = 175
$ws1.Cells.Item(26, 6).Value = 1254
$ws1.Cells.Item(28, 6).Value = 985
$ws1.Cells.Item(29, 6).Value = 333
$ws1.Cells.Item(30, 6).Value = 1931
$ws1.Cells.Item(31, 6).Value = 40
$ws1.Cells.Item(33, 6).Value = 1358

# Sheet2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value = 514
$ws2.Cells.Item(3, 7).Value = "已售罄"
$ws2.Cells.Item(4, 6).Value = 356
$ws2.Cells.Item(5, 6).Value = 10

# Sheet3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 725

# Sheet4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 725
$ws4.Cells.Item(4, 6).Value = 521
$ws4.Cells.Item(6, 6).Value = 393
$ws4.Cells.Item(7, 6).Value = 497
$ws4.Cells.Item(10, 6).Value = 356
$ws4.Cells.Item(11, 6).Value = 10
$ws4.Cells.Item(13, 6).Value = 935
$ws4.Cells.Item(16, 6).Value = 12
$ws4.Cells.Item(17, 6).Value = 975
$ws4.Cells.Item(18, 6).Value = 767
$ws4.Cells.Item(19, 6).Value = 212
$ws4.Cells.Item(27, 6).Value = 797
$ws4.Cells.Item(28, 6).Value = 262
$ws4.Cells.Item(29, 6).Value = 565
$ws4.Cells.Item(30, 6).Value = 496
$ws4.Cells.Item(31, 6).Value = 1310
$ws4.Cells.Item(33, 6).Value = 438
$ws4.Cells.Item(34, 6).Value = 1134
$ws4.Cells.Item(35, 6).Value = 2823
$ws4.Cells.Item(36, 6).Value = 1342
$ws4.Cells.Item(37, 6).Value = 670
$ws4.Cells.Item(38, 6).Value = 175
$ws4.Cells.Item(39, 6).Value = 1254
$ws4.Cells.Item(43, 6).Value = 985
$ws4.Cells.Item(44, 6).Value = 333
$ws4.Cells.Item(45, 6).Value = 1931
$ws4.Cells.Item(46, 6).Value = 40
$ws4.Cells.Item(48, 6).Value = 1358
$ws4.Cells.Item(8, 6).Value = 514
$ws4.Cells.Item(8, 7).Value = "已售罄"
$ws4.Cells.Item(9, 6).Value = 514
$ws4.Cells.Item(9, 7).Value = "已售罄"
